# Auto-generated script to apply Betfair Back/Lay odds updates (2025-12-28 workbook)
# Updates 227 numeric cells across rows 2-21 (columns F:AO) to match the new odds snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 6).Value = 1.67  # F2: 1.66 -> 1.67
$ws.Cells.Item(2, 9).Value = 6  # I2: 6.2 -> 6
$ws.Cells.Item(2, 10).Value = 4.2  # J2: 4.3 -> 4.2
$ws.Cells.Item(2, 11).Value = 4.3  # K2: 4.4 -> 4.3
$ws.Cells.Item(2, 12).Value = 1.39  # L2: 1.38 -> 1.39
$ws.Cells.Item(2, 15).Value = 1.31  # O2: 1.3 -> 1.31
$ws.Cells.Item(2, 16).Value = 2.02  # P2: 2.06 -> 2.02
$ws.Cells.Item(2, 17).Value = 1.91  # Q2: 1.9 -> 1.91
$ws.Cells.Item(2, 18).Value = 1.39  # R2: 1.4 -> 1.39
$ws.Cells.Item(2, 19).Value = 3.35  # S2: 3.3 -> 3.35
$ws.Cells.Item(2, 20).Value = 1.93  # T2: 1.92 -> 1.93
$ws.Cells.Item(2, 23).Value = 2.46  # W2: 2.48 -> 2.46
$ws.Cells.Item(2, 25).Value = 20  # Y2: 21 -> 20
$ws.Cells.Item(2, 27).Value = 170  # AA2: 160 -> 170
$ws.Cells.Item(2, 28).Value = 8.6  # AB2: 8.4 -> 8.6
$ws.Cells.Item(2, 39).Value = 130  # AM2: 120 -> 130
$ws.Cells.Item(2, 40).Value = 10  # AN2: 9.6 -> 10
$ws.Cells.Item(2, 41).Value = 100  # AO2: 95 -> 100
$ws.Cells.Item(3, 6).Value = 23  # F3: 26 -> 23
$ws.Cells.Item(3, 7).Value = 32  # G3: 36 -> 32
$ws.Cells.Item(3, 9).Value = 1.19  # I3: 1.18 -> 1.19
$ws.Cells.Item(3, 10).Value = 7.8  # J3: 8.199999999999999 -> 7.8
$ws.Cells.Item(3, 11).Value = 10.5  # K3: 11 -> 10.5
$ws.Cells.Item(3, 13).Value = 1.02  # M3: 1.01 -> 1.02
$ws.Cells.Item(3, 14).Value = 5.3  # N3: 5.2 -> 5.3
$ws.Cells.Item(3, 16).Value = 2.5  # P3: 2.46 -> 2.5
$ws.Cells.Item(3, 17).Value = 1.51  # Q3: 1.53 -> 1.51
$ws.Cells.Item(3, 18).Value = 1.6  # R3: 1.58 -> 1.6
$ws.Cells.Item(3, 19).Value = 2.26  # S3: 2.32 -> 2.26
$ws.Cells.Item(3, 20).Value = 2.5  # T3: 2.72 -> 2.5
$ws.Cells.Item(3, 22).Value = 6  # V3: 6.4 -> 6
$ws.Cells.Item(3, 28).Value = 85  # AB3: 90 -> 85
$ws.Cells.Item(3, 29).Value = 25  # AC3: 26 -> 25
$ws.Cells.Item(3, 30).Value = 16.5  # AD3: 17 -> 16.5
$ws.Cells.Item(3, 31).Value = 19  # AE3: 19.5 -> 19
$ws.Cells.Item(3, 33).Value = 120  # AG3: 990 -> 120
$ws.Cells.Item(3, 34).Value = 70  # AH3: 80 -> 70
$ws.Cells.Item(3, 35).Value = 75  # AI3: 85 -> 75
$ws.Cells.Item(4, 6).Value = 2.5  # F4: 2.48 -> 2.5
$ws.Cells.Item(4, 8).Value = 2.7  # H4: 2.58 -> 2.7
$ws.Cells.Item(4, 9).Value = 3.35  # I4: 3.55 -> 3.35
$ws.Cells.Item(4, 10).Value = 3.05  # J4: 2.8 -> 3.05
$ws.Cells.Item(4, 11).Value = 3.65  # K4: 4.8 -> 3.65
$ws.Cells.Item(4, 13).Value = 1.08  # M4: 1.07 -> 1.08
$ws.Cells.Item(4, 16).Value = 1.7  # P4: 1.69 -> 1.7
$ws.Cells.Item(4, 17).Value = 2.1  # Q4: 1.98 -> 2.1
$ws.Cells.Item(4, 20).Value = 1.81  # T4: 1.8 -> 1.81
$ws.Cells.Item(4, 22).Value = 1.43  # V4: 1.41 -> 1.43
$ws.Cells.Item(5, 6).Value = 6.6  # F5: 6.2 -> 6.6
$ws.Cells.Item(5, 7).Value = 8  # G5: 7.4 -> 8
$ws.Cells.Item(5, 8).Value = 1.49  # H5: 1.5 -> 1.49
$ws.Cells.Item(5, 12).Value = 1.33  # L5: 1.28 -> 1.33
$ws.Cells.Item(5, 16).Value = 2.16  # P5: 2.18 -> 2.16
$ws.Cells.Item(5, 18).Value = 1.46  # R5: 1.45 -> 1.46
$ws.Cells.Item(5, 23).Value = 1.14  # W5: 1.16 -> 1.14
$ws.Cells.Item(5, 25).Value = 11.5  # Y5: 12 -> 11.5
$ws.Cells.Item(5, 32).Value = 70  # AF5: 980 -> 70
$ws.Cells.Item(5, 33).Value = 1000  # AG5: 27 -> 1000
$ws.Cells.Item(5, 37).Value = 110  # AK5: 100 -> 110
$ws.Cells.Item(5, 40).Value = 130  # AN5: 120 -> 130
$ws.Cells.Item(6, 6).Value = 1.4  # F6: 1.41 -> 1.4
$ws.Cells.Item(6, 7).Value = 1.44  # G6: 1.43 -> 1.44
$ws.Cells.Item(6, 8).Value = 8.6  # H6: 8.199999999999999 -> 8.6
$ws.Cells.Item(6, 13).Value = 1.05  # M6: 1.04 -> 1.05
$ws.Cells.Item(6, 30).Value = 1000  # AD6: 40 -> 1000
$ws.Cells.Item(6, 36).Value = 13  # AJ6: 15 -> 13
$ws.Cells.Item(7, 9).Value = 5  # I7: 5.1 -> 5
$ws.Cells.Item(7, 11).Value = 3.95  # K7: 3.9 -> 3.95
$ws.Cells.Item(7, 20).Value = 1.86  # T7: 1.85 -> 1.86
$ws.Cells.Item(7, 23).Value = 2.1  # W7: 2.12 -> 2.1
$ws.Cells.Item(8, 8).Value = 9.800000000000001  # H8: 10 -> 9.800000000000001
$ws.Cells.Item(8, 19).Value = 3.35  # S8: 3.4 -> 3.35
$ws.Cells.Item(8, 25).Value = 28  # Y8: 27 -> 28
$ws.Cells.Item(8, 31).Value = 190  # AE8: 180 -> 190
$ws.Cells.Item(8, 32).Value = 7.4  # AF8: 7.6 -> 7.4
$ws.Cells.Item(9, 7).Value = 2.48  # G9: 2.46 -> 2.48
$ws.Cells.Item(9, 8).Value = 3.25  # H9: 3.3 -> 3.25
$ws.Cells.Item(9, 9).Value = 3.65  # I9: 3.6 -> 3.65
$ws.Cells.Item(9, 22).Value = 1.39  # V9: 1.38 -> 1.39
$ws.Cells.Item(9, 23).Value = 1.68  # W9: 1.69 -> 1.68
$ws.Cells.Item(9, 24).Value = 13  # X9: 14.5 -> 13
$ws.Cells.Item(10, 8).Value = 2.08  # H10: 2.1 -> 2.08
$ws.Cells.Item(10, 9).Value = 2.2  # I10: 2.22 -> 2.2
$ws.Cells.Item(10, 10).Value = 3.3  # J10: 3.35 -> 3.3
$ws.Cells.Item(10, 22).Value = 1.83  # V10: 1.81 -> 1.83
$ws.Cells.Item(11, 6).Value = 7.2  # F11: 7 -> 7.2
$ws.Cells.Item(11, 7).Value = 7.4  # G11: 7.2 -> 7.4
$ws.Cells.Item(11, 10).Value = 4.1  # J11: 4.2 -> 4.1
$ws.Cells.Item(11, 11).Value = 4.2  # K11: 4.3 -> 4.2
$ws.Cells.Item(11, 15).Value = 1.38  # O11: 1.39 -> 1.38
$ws.Cells.Item(11, 22).Value = 2.6  # V11: 2.62 -> 2.6
$ws.Cells.Item(11, 27).Value = 14.5  # AA11: 15 -> 14.5
$ws.Cells.Item(11, 33).Value = 27  # AG11: 28 -> 27
$ws.Cells.Item(11, 35).Value = 46  # AI11: 44 -> 46
$ws.Cells.Item(12, 10).Value = 3.2  # J12: 3.25 -> 3.2
$ws.Cells.Item(12, 11).Value = 3.25  # K12: 3.3 -> 3.25
$ws.Cells.Item(12, 19).Value = 4.9  # S12: 5.1 -> 4.9
$ws.Cells.Item(12, 20).Value = 2.08  # T12: 2.06 -> 2.08
$ws.Cells.Item(12, 24).Value = 9.4  # X12: 9.6 -> 9.4
$ws.Cells.Item(12, 27).Value = 46  # AA12: 48 -> 46
$ws.Cells.Item(12, 28).Value = 8.800000000000001  # AB12: 9 -> 8.800000000000001
$ws.Cells.Item(12, 30).Value = 13.5  # AD12: 14 -> 13.5
$ws.Cells.Item(12, 31).Value = 36  # AE12: 38 -> 36
$ws.Cells.Item(12, 37).Value = 36  # AK12: 38 -> 36
$ws.Cells.Item(12, 38).Value = 60  # AL12: 65 -> 60
$ws.Cells.Item(12, 40).Value = 42  # AN12: 44 -> 42
$ws.Cells.Item(12, 41).Value = 42  # AO12: 44 -> 42
$ws.Cells.Item(13, 9).Value = 2.06  # I13: 2.08 -> 2.06
$ws.Cells.Item(13, 14).Value = 3.05  # N13: 3.1 -> 3.05
$ws.Cells.Item(13, 20).Value = 1.99  # T13: 2 -> 1.99
$ws.Cells.Item(13, 22).Value = 1.94  # V13: 1.93 -> 1.94
$ws.Cells.Item(14, 6).Value = 3.4  # F14: 3.45 -> 3.4
$ws.Cells.Item(14, 7).Value = 3.5  # G14: 3.55 -> 3.5
$ws.Cells.Item(14, 8).Value = 2.48  # H14: 2.44 -> 2.48
$ws.Cells.Item(14, 9).Value = 2.54  # I14: 2.52 -> 2.54
$ws.Cells.Item(14, 10).Value = 3.15  # J14: 3.2 -> 3.15
$ws.Cells.Item(14, 11).Value = 3.2  # K14: 3.25 -> 3.2
$ws.Cells.Item(14, 17).Value = 2.56  # Q14: 2.58 -> 2.56
$ws.Cells.Item(14, 22).Value = 1.64  # V14: 1.66 -> 1.64
$ws.Cells.Item(14, 23).Value = 1.4  # W14: 1.39 -> 1.4
$ws.Cells.Item(14, 24).Value = 8.6  # X14: 9 -> 8.6
$ws.Cells.Item(14, 26).Value = 16  # Z14: 17.5 -> 16
$ws.Cells.Item(14, 27).Value = 38  # AA14: 980 -> 38
$ws.Cells.Item(14, 30).Value = 12.5  # AD14: 13 -> 12.5
$ws.Cells.Item(14, 31).Value = 34  # AE14: 980 -> 34
$ws.Cells.Item(14, 32).Value = 24  # AF14: 25 -> 24
$ws.Cells.Item(14, 33).Value = 15.5  # AG14: 16.5 -> 15.5
$ws.Cells.Item(14, 36).Value = 75  # AJ14: 80 -> 75
$ws.Cells.Item(14, 37).Value = 55  # AK14: 60 -> 55
$ws.Cells.Item(14, 39).Value = 180  # AM14: 200 -> 180
$ws.Cells.Item(14, 40).Value = 75  # AN14: 80 -> 75
$ws.Cells.Item(14, 41).Value = 36  # AO14: 980 -> 36
$ws.Cells.Item(15, 6).Value = 1.86  # F15: 1.98 -> 1.86
$ws.Cells.Item(15, 7).Value = 1.99  # G15: 2.12 -> 1.99
$ws.Cells.Item(15, 8).Value = 4  # H15: 3.65 -> 4
$ws.Cells.Item(15, 9).Value = 4.7  # I15: 4.2 -> 4.7
$ws.Cells.Item(15, 10).Value = 3.8  # J15: 3.65 -> 3.8
$ws.Cells.Item(15, 17).Value = 1.73  # Q15: 1.74 -> 1.73
$ws.Cells.Item(15, 20).Value = 1.69  # T15: 1.67 -> 1.69
$ws.Cells.Item(15, 21).Value = 2.2  # U15: 2.24 -> 2.2
$ws.Cells.Item(15, 22).Value = 1.27  # V15: 1.32 -> 1.27
$ws.Cells.Item(15, 23).Value = 2  # W15: 1.89 -> 2
$ws.Cells.Item(15, 26).Value = 42  # Z15: 34 -> 42
$ws.Cells.Item(15, 28).Value = 11  # AB15: 13 -> 11
$ws.Cells.Item(15, 30).Value = 18  # AD15: 19.5 -> 18
$ws.Cells.Item(15, 32).Value = 16  # AF15: 17 -> 16
$ws.Cells.Item(15, 33).Value = 12.5  # AG15: 13 -> 12.5
$ws.Cells.Item(15, 36).Value = 27  # AJ15: 30 -> 27
$ws.Cells.Item(15, 37).Value = 23  # AK15: 25 -> 23
$ws.Cells.Item(15, 38).Value = 38  # AL15: 40 -> 38
$ws.Cells.Item(15, 40).Value = 14  # AN15: 15.5 -> 14
$ws.Cells.Item(16, 6).Value = 2.3  # F16: 2.28 -> 2.3
$ws.Cells.Item(16, 8).Value = 3.55  # H16: 3.6 -> 3.55
$ws.Cells.Item(16, 9).Value = 3.6  # I16: 3.65 -> 3.6
$ws.Cells.Item(16, 15).Value = 1.35  # O16: 1.34 -> 1.35
$ws.Cells.Item(16, 16).Value = 1.93  # P16: 1.91 -> 1.93
$ws.Cells.Item(16, 17).Value = 2.04  # Q16: 2.06 -> 2.04
$ws.Cells.Item(16, 20).Value = 1.84  # T16: 1.83 -> 1.84
$ws.Cells.Item(16, 23).Value = 1.75  # W16: 1.76 -> 1.75
$ws.Cells.Item(16, 26).Value = 23  # Z16: 24 -> 23
$ws.Cells.Item(16, 28).Value = 9.800000000000001  # AB16: 9.6 -> 9.800000000000001
$ws.Cells.Item(16, 31).Value = 38  # AE16: 42 -> 38
$ws.Cells.Item(16, 38).Value = 38  # AL16: 40 -> 38
$ws.Cells.Item(16, 40).Value = 18.5  # AN16: 19 -> 18.5
$ws.Cells.Item(16, 41).Value = 40  # AO16: 42 -> 40
$ws.Cells.Item(17, 7).Value = 1.84  # G17: 1.83 -> 1.84
$ws.Cells.Item(17, 10).Value = 3.7  # J17: 3.75 -> 3.7
$ws.Cells.Item(17, 21).Value = 2  # U17: 1.99 -> 2
$ws.Cells.Item(17, 23).Value = 2.18  # W17: 2.2 -> 2.18
$ws.Cells.Item(17, 25).Value = 17  # Y17: 17.5 -> 17
$ws.Cells.Item(17, 27).Value = 130  # AA17: 150 -> 130
$ws.Cells.Item(17, 30).Value = 20  # AD17: 21 -> 20
$ws.Cells.Item(17, 35).Value = 80  # AI17: 85 -> 80
$ws.Cells.Item(17, 40).Value = 13  # AN17: 12.5 -> 13
$ws.Cells.Item(18, 7).Value = 3.85  # G18: 4 -> 3.85
$ws.Cells.Item(18, 8).Value = 2.2  # H18: 2.16 -> 2.2
$ws.Cells.Item(18, 10).Value = 3.5  # J18: 3.45 -> 3.5
$ws.Cells.Item(18, 14).Value = 3.4  # N18: 3.35 -> 3.4
$ws.Cells.Item(18, 17).Value = 2.08  # Q18: 2.12 -> 2.08
$ws.Cells.Item(18, 22).Value = 1.81  # V18: 1.83 -> 1.81
$ws.Cells.Item(18, 23).Value = 1.35  # W18: 1.34 -> 1.35
$ws.Cells.Item(18, 24).Value = 12.5  # X18: 12 -> 12.5
$ws.Cells.Item(18, 30).Value = 10.5  # AD18: 11 -> 10.5
$ws.Cells.Item(18, 33).Value = 15.5  # AG18: 16 -> 15.5
$ws.Cells.Item(18, 34).Value = 19  # AH18: 19.5 -> 19
$ws.Cells.Item(18, 40).Value = 55  # AN18: 1000 -> 55
$ws.Cells.Item(18, 41).Value = 20  # AO18: 19.5 -> 20
$ws.Cells.Item(19, 8).Value = 2.08  # H19: 2.1 -> 2.08
$ws.Cells.Item(19, 9).Value = 2.1  # I19: 2.12 -> 2.1
$ws.Cells.Item(19, 14).Value = 4.3  # N19: 4.5 -> 4.3
$ws.Cells.Item(19, 15).Value = 1.28  # O19: 1.27 -> 1.28
$ws.Cells.Item(19, 16).Value = 2.14  # P19: 2.16 -> 2.14
$ws.Cells.Item(19, 20).Value = 1.75  # T19: 1.73 -> 1.75
$ws.Cells.Item(19, 21).Value = 2.28  # U19: 2.32 -> 2.28
$ws.Cells.Item(19, 22).Value = 1.9  # V19: 1.89 -> 1.9
$ws.Cells.Item(19, 27).Value = 24  # AA19: 25 -> 24
$ws.Cells.Item(19, 28).Value = 16.5  # AB19: 16 -> 16.5
$ws.Cells.Item(19, 39).Value = 80  # AM19: 75 -> 80
$ws.Cells.Item(20, 7).Value = 1.16  # G20: 1.17 -> 1.16
$ws.Cells.Item(20, 9).Value = 27  # I20: 26 -> 27
$ws.Cells.Item(20, 10).Value = 10  # J20: 9.4 -> 10
$ws.Cells.Item(20, 11).Value = 11.5  # K20: 10.5 -> 11.5
$ws.Cells.Item(20, 14).Value = 7.4  # N20: 7 -> 7.4
$ws.Cells.Item(20, 16).Value = 3.15  # P20: 3.1 -> 3.15
$ws.Cells.Item(20, 19).Value = 2.06  # S20: 2.04 -> 2.06
$ws.Cells.Item(20, 20).Value = 2.32  # T20: 2.28 -> 2.32
$ws.Cells.Item(20, 21).Value = 1.69  # U20: 1.68 -> 1.69
$ws.Cells.Item(20, 22).Value = 1.03  # V20: 1.04 -> 1.03
$ws.Cells.Item(20, 23).Value = 7.2  # W20: 6.8 -> 7.2
$ws.Cells.Item(20, 24).Value = 55  # X20: 990 -> 55
$ws.Cells.Item(20, 25).Value = 95  # Y20: 990 -> 95
$ws.Cells.Item(20, 26).Value = 310  # Z20: 280 -> 310
$ws.Cells.Item(20, 28).Value = 12  # AB20: 12.5 -> 12
$ws.Cells.Item(20, 29).Value = 27  # AC20: 23 -> 27
$ws.Cells.Item(20, 30).Value = 100  # AD20: 990 -> 100
$ws.Cells.Item(20, 31).Value = 520  # AE20: 460 -> 520
$ws.Cells.Item(20, 32).Value = 8.800000000000001  # AF20: 9.199999999999999 -> 8.800000000000001
$ws.Cells.Item(20, 33).Value = 13.5  # AG20: 14 -> 13.5
$ws.Cells.Item(20, 35).Value = 340  # AI20: 300 -> 340
$ws.Cells.Item(20, 37).Value = 14.5  # AK20: 15.5 -> 14.5
$ws.Cells.Item(20, 38).Value = 55  # AL20: 980 -> 55
$ws.Cells.Item(20, 39).Value = 340  # AM20: 280 -> 340
$ws.Cells.Item(21, 10).Value = 4.1  # J21: 3.85 -> 4.1
$ws.Cells.Item(21, 14).Value = 2.84  # N21: 2.82 -> 2.84
$ws.Cells.Item(21, 16).Value = 2.04  # P21: 1.92 -> 2.04
$ws.Cells.Item(21, 18).Value = 1.3  # R21: 1.34 -> 1.3
$ws.Cells.Item(21, 19).Value = 2.28  # S21: 2.52 -> 2.28

Write-Host "Applied 227 cell updates"
